$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D13").Select()
